$d = $word.ActiveDocument

$newText = "You are participating in a global campaign to observe and record the faintest stars visible as a means of measuring light pollution in a given location. By locating and observing the constellation Leo constellation in the night sky and comparing it to stellar charts, people from around the world will learn how the lights in their community contribute to light pollution. Your contributions to the online database will document the visible nighttime sky."

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*You are participating*By locating and observing the constellation Perseus*") {
        $r = $p.Range
        # Exclude the trailing paragraph mark so we only replace the run content.
        $r.End = $r.End - 1
        $r.Delete()
        $r.InsertAfter($newText)
    }
}
